# New trade row appended to the sheet (row 6).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give G6 the same date/time number format already used by G2:G5 before
# writing its value, so it reuses the existing style instead of minting a
# new one.
$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A6").Value = 9899.1200000000008
$ws.Range("B6").Value = 9975.93
$ws.Range("C6").Value = 79.650000000000006
$ws.Range("D6").Value = 79.040000000000006
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = -0.77
$ws.Range("G6").Value = 42612.67423611111
$ws.Range("H6").Value = $false
